# Update "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" worksheets, reflecting the regenerated output
# published to gh-pages (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Worksheet "展览"
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 304
$wsExpo.Range("F4").Value = 3694
$wsExpo.Range("F5").Value = 2254
$wsExpo.Range("F8").Value = 9
$wsExpo.Range("F12").Value = 1367
$wsExpo.Range("F13").Value = 240
$wsExpo.Range("F14").Value = 2115
$wsExpo.Range("F15").Value = 154

# Worksheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 304
$wsAll.Range("F4").Value = 3694
$wsAll.Range("F5").Value = 2254
$wsAll.Range("F8").Value = 9
$wsAll.Range("F15").Value = 1367
$wsAll.Range("F16").Value = 240
$wsAll.Range("F17").Value = 2115
$wsAll.Range("F18").Value = 154
